# publish 0.0.5 sheet2excel add para:output_path
#
# before.xlsx has a single sheet ("Sheet1") with a two-row "name" / "严健"
# column A.  The target workbook adds two more sheets ("Sheet2", "Sheet3"),
# each holding a single extra data point, and leaves Sheet3 as the
# selected/active sheet.

$wb = $excel.ActiveWorkbook

# Sheet1 already exists and already holds the "name" / "严健" values with
# the original cell styling - leave it untouched.
$ws1 = $wb.Worksheets.Item(1)

# Add Sheet2 right after Sheet1, then Sheet3 right after Sheet2, so the
# final tab order is Sheet1, Sheet2, Sheet3.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws3 = $wb.Worksheets.Add($null, $ws2)

# Sheet2: A1 = 富士达
$ws2.Range("A1").Value = "富士达"

# Sheet2 page setup: A4 paper, portrait orientation.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Sheet3: C10 = 发大水
$ws3.Range("C10").Value = "发大水"

# Leave the cursor/selection where the author left it, and make Sheet3 the
# active (selected) sheet, matching the saved workbook view state.
[void]$ws2.Range("D16").Select()
[void]$ws3.Range("C10").Select()
